# Added passwords for every account
#
# UserList sheet gains an "active" flag (column G) and a password-hash
# column (column H) for every user row. The previously-empty row 5 is
# filled in with account "e" (which used to sit at row 6), so accounts
# "f"/"g"/"banana" each shift up one row (6->5 is "e" itself, 7->6,
# 8->7, 9->8), and every row 1-8 ends up with a flag + hash in G/H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserList")

# Row 1 (cmoticska) already had G1/H1 - just flip the flag to active.
$ws.Cells.Item(1, 7).Value = 1

# Row 2 (jedwards) gains the active flag.
$ws.Cells.Item(2, 7).Value = 1

# Row 3 (mcolumbo): flag it active and replace its password hash.
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = "6cf615d5bcaac778352a8f1f3360d23f02f34ec182e259897fd6ce485d7870d4"

# Row 4 (arivera): new active flag + new password hash.
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = "5906ac361a137e2d286465cd6588ebb5ac3f5ae955001100bc41577c3d751764"

# Row 5 was an empty gap; fill it with account "e" (previously at row 6),
# plus its flag + hash.
$ws.Cells.Item(5, 1).Value = "e"
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = "b97873a40f73abedd8d685a7cd5e5f85e4a9cfb83eac26886640a0813850122b"

# Account "f" moves up from row 7 to row 6, plus flag + hash.
$ws.Cells.Item(6, 1).Value = "f"
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = "8b2c86ea9cf2ea4eb517fd1e06b74f399e7fec0fef92e3b482a6cf2e2b092023"

# Account "g" moves up from row 8 to row 7, plus flag + hash.
$ws.Cells.Item(7, 1).Value = "g"
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = "598a1a400c1dfdf36974e69d7e1bc98593f2e15015eed8e9b7e47a83b31693d5"

# Account "banana" moves up from row 9 to row 8, plus flag + hash
# (replacing the stray hash that used to live at H9).
$ws.Cells.Item(8, 1).Value = "banana"
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = "5860836e8f13fc9837539a597d4086bfc0299e54ad92148d54538b5c3feefb7c"

# Row 9 no longer holds any data now that "banana" moved up to row 8.
$ws.Rows.Item(9).ClearContents()

# Update the active selection to match the author's final cursor position.
$ws.Range("G4").Select() | Out-Null
